$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial number that was bumped from
# 45186 (2023-09-17) to 45188 (2023-09-19) for every data row (rows 2-257).
for ($row = 2; $row -le 257; $row++) {
    $ws.Cells.Item($row, 3).Value = 45188
}
